# Update cryptocurrency price/volume table to latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold numeric-looking text (prices/percentages) that must stay
# plain text (inlineStr), exactly like the rest of the sheet, instead of
# being auto-coerced to numbers by Excel. Force text format for the whole
# touched range first, write the values, then drop the temporary format so
# the cells end up back on the sheet's default (unstyled) text cells.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.503.96'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.726.76'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '245.16'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '0.4804'
$ws.Range('E7').Value = '  +1.56%  '
$ws.Range('D8').Value = '0.2668'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').Value = '0.06218'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').Value = '1.725.75'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('D11').Value = '0.07147'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').Value = '15.65'
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').Value = '0.6162'
$ws.Range('E13').Value = '  +3.74%  '
$ws.Range('D14').Value = '4.520'
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').Value = '77.09'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = '26.505.13'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '1.0000'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '0.000006936'
$ws.Range('E19').Value = '  +1.85%  '
$ws.Range('D20').Value = '11.66'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').Value = '1.947.37'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').Value = '4.521'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').Value = '8.946'
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('D24').Value = '5.277'
$ws.Range('E24').Value = '  -1.21%  '
$ws.Range('D25').Value = '136.38'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').Value = '15.33'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').Value = '1.792'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').Value = '1.405'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').Value = '106.75'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = '3.974'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').Value = '0.08021'
$ws.Range('E31').Value = '  +3.36%  '
$ws.Range('D32').Value = '3.703'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').Value = '0.9996'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.616'
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.6372'
$ws.Range('E36').Value = '  +2.53%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '0.9905'
$ws.Range('E37').Value = '  +1.24%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '0.9288'
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '2.089'
$ws.Range('E39').Value = '  +9.16%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.417'
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').Value = '104.92'
$ws.Range('E41').Value = '  -9.60%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.006'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.01502'
$ws.Range('E43').Value = '  +1.61%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.630'
$ws.Range('E44').Value = '  +4.93%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.3903'
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.912'
$ws.Range('E46').Value = '  +10.09%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1184'
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.05331'
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '30.87'
$ws.Range('E49').Value = '  +0.94%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.854'
$ws.Range('E50').Value = '  +2.46%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.269'
$ws.Range('E51').Value = '  +3.93%  '

$numRange.Style = "Normal"
